$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: split the run "-in each of the 6 demographics" into two runs:
#   "-in each of the"  +  " 6 demographics"
# ---------------------------------------------------------------------------

$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "-in each of the 6 demographics") {
        $target1 = $p
        break
    }
}

if ($target1 -ne $null) {
    $xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>-in each of the</w:t></w:r><w:r><w:t xml:space="preserve"> 6 demographics</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $target1.Range.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Edit 2: append two new bold paragraphs at the very end of the document,
# right after the last (empty) paragraph and before the section break.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>* Provide guidance on how to exclude the &#8216;No Answer&#8217; counts</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">* Enable deeper cross-tabulation.  E.g. Of those who identified as black, investigate why they are or aren&#8217;t interested in HMS Promotion; are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>the majority of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> them long-time Instructors?</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$lastPara.Range.InsertXML($xml2)
